$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Apollo 11")

$ws.Range("B3").Value = 44.8
$ws.Range("B4").Value = 44.4
$ws.Range("B5").Value = 78.400000000000006
$ws.Range("B6").Value = 78.3

$ws.Range("B10").Value = 15712
$ws.Range("B11").Value = 25091

$ws.Range("B12").Value = 110
$ws.Range("B13").Value = 110
$ws.Range("B14").Value = 110
$ws.Range("B15").Value = 110
$ws.Range("B16").Value = 225
$ws.Range("B17").Value = 225
$ws.Range("B18").Value = 225
$ws.Range("B19").Value = 225

$ws.Range("B26").Value = 6975
$ws.Range("B27").Value = 11209

$ws.Range("B31").Value = 2020
$ws.Range("B32").Value = 3218
$ws.Range("B33").Value = 108
$ws.Range("B34").Value = 108
$ws.Range("B35").Value = 209
$ws.Range("B36").Value = 209

$ws.Range("A19").Select()
$ws.Range("B40").Select()
